$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C6").Value = 0
$ws.Range("H58").Formula = "=H53+H56+H57"
$ws.Range("H60").Formula = "=H26+H34+H58"

$excel.ActiveWindow.ScrollRow = 19
$ws.Range("G24").Select()
